$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Mmp9"
$ws.Range("C2").Value2 = "Lrp1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.04273
$ws.Range("H2").Value2 = 0.12819
$ws.Range("I2").Value2 = 0.7009974407769539
$ws.Range("J2").Value2 = 0.7009974407769538
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 3.456265333333333
$ws.Range("N2").Value2 = 10.368796
$ws.Range("O2").Value2 = 0.009841535807677501
$ws.Range("P2").Value2 = 0.0098415358076775
$ws.Range("Q2").Value2 = 0.1476862176933333
$ws.Range("R2").Value2 = 1.32917595924
$ws.Range("S2").Value2 = 0.006898891414496681
$ws.Range("T2").Value2 = 0.006898891414496678

# Row 3
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Mmp9"
$ws.Range("C3").Value2 = "Lrp1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.04273
$ws.Range("H3").Value2 = 0.12819
$ws.Range("I3").Value2 = 0.7009974407769539
$ws.Range("J3").Value2 = 0.7009974407769538
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 301.6001486666667
$ws.Range("N3").Value2 = 904.800446
$ws.Range("O3").Value2 = 0.8587907398420774
$ws.Range("P3").Value2 = 0.8587907398420773
$ws.Range("Q3").Value2 = 12.88737435252667
$ws.Range("R3").Value2 = 115.98636917274
$ws.Range("S3").Value2 = 0.6020101107922431
$ws.Range("T3").Value2 = 0.602010110792243

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Mmp9"
$ws.Range("C4").Value2 = "Lrp1"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.04273
$ws.Range("H4").Value2 = 0.12819
$ws.Range("I4").Value2 = 0.7009974407769539
$ws.Range("J4").Value2 = 0.7009974407769538
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 46.13524966666667
$ws.Range("N4").Value2 = 138.405749
$ws.Range("O4").Value2 = 0.1313677243502452
$ws.Range("P4").Value2 = 0.1313677243502452
$ws.Range("Q4").Value2 = 1.971359218256667
$ws.Range("R4").Value2 = 17.74223296431
$ws.Range("S4").Value2 = 0.0920884385702142
$ws.Range("T4").Value2 = 0.09208843857021418

# Row 5
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Mmp9"
$ws.Range("C5").Value2 = "Lrp1"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.018226
$ws.Range("H5").Value2 = 0.054678
$ws.Range("I5").Value2 = 0.2990025592230461
$ws.Range("J5").Value2 = 0.2990025592230461
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 3.456265333333333
$ws.Range("N5").Value2 = 10.368796
$ws.Range("O5").Value2 = 0.009841535807677501
$ws.Range("P5").Value2 = 0.0098415358076775
$ws.Range("Q5").Value2 = 0.06299389196533332
$ws.Range("R5").Value2 = 0.5669450276879999
$ws.Range("S5").Value2 = 0.002942644393180821
$ws.Range("T5").Value2 = 0.002942644393180821

# Row 6
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Mmp9"
$ws.Range("C6").Value2 = "Lrp1"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.018226
$ws.Range("H6").Value2 = 0.054678
$ws.Range("I6").Value2 = 0.2990025592230461
$ws.Range("J6").Value2 = 0.2990025592230461
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 301.6001486666667
$ws.Range("N6").Value2 = 904.800446
$ws.Range("O6").Value2 = 0.8587907398420774
$ws.Range("P6").Value2 = 0.8587907398420773
$ws.Range("Q6").Value2 = 5.496964309598666
$ws.Range("R6").Value2 = 49.47267878638799
$ws.Range("S6").Value2 = 0.2567806290498343
$ws.Range("T6").Value2 = 0.2567806290498343

# Row 7
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Mmp9"
$ws.Range("C7").Value2 = "Lrp1"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.018226
$ws.Range("H7").Value2 = 0.054678
$ws.Range("I7").Value2 = 0.2990025592230461
$ws.Range("J7").Value2 = 0.2990025592230461
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 46.13524966666667
$ws.Range("N7").Value2 = 138.405749
$ws.Range("O7").Value2 = 0.1313677243502452
$ws.Range("P7").Value2 = 0.1313677243502452
$ws.Range("Q7").Value2 = 0.8408610604246668
$ws.Range("R7").Value2 = 7.567749543822001
$ws.Range("S7").Value2 = 0.03927928578003099
$ws.Range("T7").Value2 = 0.03927928578003098

# Remove old trailing rows (8-10) that no longer exist in the refreshed output
$ws.Rows("8:10").Delete()

Write-Host "done"